$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.2589695869521966
$ws.Range("D2").Value = 0.7980663936084045

$ws.Range("C3").Value = 2.837697657568838
$ws.Range("D3").Value = 0.009577714291172823

$ws.Range("C4").Value = 0.7282625919490749
$ws.Range("D4").Value = 0.4741338626895228

$ws.Range("C5").Value = -0.05209187286232511
$ws.Range("D5").Value = 0.9589255961969929

$ws.Range("C6").Value = 1.454963415918597
$ws.Range("D6").Value = 0.1598003291846837

$ws.Range("C7").Value = 0.4239998253858546
$ws.Range("D7").Value = 0.6756830419372992

$ws.Range("C8").Value = -0.2345623162375872
$ws.Range("D8").Value = 0.8167177561495946

$ws.Range("C9").Value = -1.167079326485868
$ws.Range("D9").Value = 0.2556730650255756

$ws.Range("C10").Value = -2.042012647552816
$ws.Range("D10").Value = 0.05331902861636562

$ws.Range("C11").Value = -0.6775430290495614
$ws.Range("D11").Value = 0.5051250668620766
